$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 360
$ws.Range("F2").Value = 45992
$ws.Range("G2").Value = 30865
$ws.Range("H2").Value = 46055

# Row 3
$ws.Range("E3").Value = 30803
$ws.Range("F3").Value = 45992

# Row 4
$ws.Range("E4").Value = 30803
$ws.Range("F4").Value = 45992
$ws.Range("G4").Value = 30865
$ws.Range("H4").Value = 46055

# Row 5
$ws.Range("E5").Value = 30803
$ws.Range("F5").Value = 45992
$ws.Range("G5").Value = 30865
$ws.Range("H5").Value = 46055

# Row 6
$ws.Range("G6").Value = 30865
$ws.Range("H6").Value = 46055

# Row 7
$ws.Range("E7").Value = 30773
$ws.Range("F7").Value = 45962
$ws.Range("G7").Value = 30865
$ws.Range("H7").Value = 46055

# Row 8
$ws.Range("D8").Value = 436
$ws.Range("H8").Value = 46055

# Row 9
$ws.Range("G9").Value = 30865
$ws.Range("H9").Value = 46055

# Row 10
$ws.Range("E10").Value = 30803
$ws.Range("F10").Value = 45992
$ws.Range("G10").Value = 30865
$ws.Range("H10").Value = 46055

# Row 11
$ws.Range("E11").Value = 30773
$ws.Range("F11").Value = 45962
$ws.Range("G11").Value = 30865
$ws.Range("H11").Value = 46055

# Row 12
$ws.Range("C12").Value = 397
$ws.Range("D12").Value = 378
$ws.Range("F12").Value = 45992
$ws.Range("H12").Value = 46055

# Row 13
$ws.Range("C13").Value = 493
$ws.Range("F13").Value = 45992
$ws.Range("G13").Value = 30865
$ws.Range("H13").Value = 46055

# Row 14
$ws.Range("C14").Value = 448
$ws.Range("D14").Value = 422
$ws.Range("F14").Value = 45962
$ws.Range("H14").Value = 46055

# Row 15
$ws.Range("C15").Value = 409
$ws.Range("F15").Value = 45962
$ws.Range("G15").Value = 30834
$ws.Range("H15").Value = 46055

# Row 16
$ws.Range("D16").Value = 436
$ws.Range("H16").Value = 46055

# Row 17
$ws.Range("C17").Value = 396
$ws.Range("D17").Value = 420
$ws.Range("F17").Value = 45992
$ws.Range("H17").Value = 46055

# Row 18
$ws.Range("D18").Value = 286
$ws.Range("E18").Value = 30803
$ws.Range("F18").Value = 45992
$ws.Range("H18").Value = 46055

# Row 19
$ws.Range("D19").Value = 424
$ws.Range("E19").Value = 30803
$ws.Range("F19").Value = 45992
$ws.Range("H19").Value = 46055

# Row 20
$ws.Range("E20").Value = 28460
$ws.Range("F20").Value = 45992
$ws.Range("G20").Value = 30865
$ws.Range("H20").Value = 46055

# Row 21
$ws.Range("C21").Value = 336
$ws.Range("F21").Value = 45992
$ws.Range("G21").Value = 30865
$ws.Range("H21").Value = 46055

# Row 22
$ws.Range("C22").Value = 349
$ws.Range("D22").Value = 393
$ws.Range("F22").Value = 45992
$ws.Range("H22").Value = 46055

# Row 23
$ws.Range("C23").Value = 286
$ws.Range("F23").Value = 45962

# Row 24
$ws.Range("D24").Value = 436
$ws.Range("H24").Value = 46055

# Row 25
$ws.Range("D25").Value = 341
$ws.Range("E25").Value = 30742
$ws.Range("F25").Value = 45931
$ws.Range("H25").Value = 46055

# Row 26
$ws.Range("D26").Value = 339
$ws.Range("H26").Value = 46055

# Row 27
$ws.Range("E27").Value = 30803
$ws.Range("F27").Value = 45992
$ws.Range("G27").Value = 30865
$ws.Range("H27").Value = 46055

# Row 28
$ws.Range("D28").Value = 406
$ws.Range("E28").Value = 30773
$ws.Range("F28").Value = 45962
$ws.Range("H28").Value = 46055

# Row 29
$ws.Range("C29").Value = 287
$ws.Range("D29").Value = 262
$ws.Range("F29").Value = 45962
$ws.Range("H29").Value = 46055

# Row 30
$ws.Range("D30").Value = 243
$ws.Range("H30").Value = 46055

# Row 31
$ws.Range("C31").Value = 420
$ws.Range("F31").Value = 45992
$ws.Range("G31").Value = 30865
$ws.Range("H31").Value = 46055

# Row 32
$ws.Range("E32").Value = 30803
$ws.Range("F32").Value = 45992
$ws.Range("G32").Value = 30865
$ws.Range("H32").Value = 46055

# Row 34
$ws.Range("D34").Value = 343
$ws.Range("H34").Value = 46055

# Row 35
$ws.Range("C35").Value = 430
$ws.Range("D35").Value = 343
$ws.Range("F35").Value = 45962
$ws.Range("H35").Value = 46055

# Row 36
$ws.Range("D36").Value = 436
$ws.Range("E36").Value = 30803
$ws.Range("F36").Value = 45992
$ws.Range("H36").Value = 46055

# Row 37
$ws.Range("D37").Value = 343
$ws.Range("H37").Value = 46055

# Row 38
$ws.Range("C38").Value = 385
$ws.Range("D38").Value = 393
$ws.Range("F38").Value = 45992
$ws.Range("H38").Value = 46055

# Row 39
$ws.Range("C39").Value = 253
$ws.Range("D39").Value = 249
$ws.Range("F39").Value = 45992
$ws.Range("H39").Value = 46055

# Row 40
$ws.Range("C40").Value = 310
$ws.Range("D40").Value = 340
$ws.Range("F40").Value = 45992
$ws.Range("H40").Value = 46055

# Row 41
$ws.Range("D41").Value = 262
$ws.Range("H41").Value = 46055

# Row 42
$ws.Range("C42").Value = 262
$ws.Range("D42").Value = 244
$ws.Range("F42").Value = 45962
$ws.Range("H42").Value = 46055

# Row 43
$ws.Range("D43").Value = 343
$ws.Range("H43").Value = 46055

# Row 44
$ws.Range("C44").Value = 431
$ws.Range("D44").Value = 330
$ws.Range("F44").Value = 45962
$ws.Range("H44").Value = 46055

# Row 45
$ws.Range("D45").Value = 343
$ws.Range("H45").Value = 46055

# Row 46
$ws.Range("D46").Value = 324
$ws.Range("H46").Value = 46055

# Row 47
$ws.Range("C47").Value = 361
$ws.Range("D47").Value = 286
$ws.Range("F47").Value = 45992
$ws.Range("H47").Value = 46055

# Row 48
$ws.Range("D48").Value = 341
$ws.Range("H48").Value = 46055

# Row 49
$ws.Range("D49").Value = 339
$ws.Range("H49").Value = 46055

# Row 50
$ws.Range("C50").Value = 384
$ws.Range("D50").Value = 264
$ws.Range("F50").Value = 45992
$ws.Range("G50").Value = 36892
$ws.Range("G50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H50").Value = 46055
$ws.Range("H50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I50").Value = "FX_limited"
$ws.Range("J50").Value = "M2_limited"
$ws.Range("K50").Value = "USDKZT"

# Row 51
$ws.Range("D51").Value = 343
$ws.Range("E51").Value = 30317
$ws.Range("H51").Value = 46055

# Row 52
$ws.Range("D52").Value = 341
$ws.Range("H52").Value = 46055
